# Recompute plate-coefficient variances (and refresh the associated
# AR/DEC coefficient columns from the updated covariance estimate).
# This mirrors the backend change that fixes the variance calculation
# for the "LIGHT_2025-02-26_20-10-15_L_1x1_-10.00_3.00s_0000" plate-coefs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (X3)
$ws.Range("C2").Value = [double]"-3.446445297159059e-15"
$ws.Range("D2").Value = [double]"5.181852886273114e-14"
$ws.Range("E2").Value = [double]"-6.015181792529995e-17"
$ws.Range("F2").Value = [double]"9.044039421943712e-16"
$ws.Range("G2").Value = [double]"2.236793990413977e-29"
$ws.Range("H2").Value = [double]"6.813664139543215e-33"
$ws.Range("I2").Value = [double]"5.63655342688161e-29"
$ws.Range("J2").Value = [double]"1.716992361386606e-32"

# Row 3 (X2Y)
$ws.Range("C3").Value = [double]"-6.833333893500072e-14"
$ws.Range("D3").Value = [double]"1.457903600507377e-14"
$ws.Range("E3").Value = [double]"-1.192641753296998e-15"
$ws.Range("F3").Value = [double]"2.544521800553381e-16"
$ws.Range("G3").Value = [double]"4.445306619349285e-29"
$ws.Range("H3").Value = [double]"1.354117832546956e-32"
$ws.Range("I3").Value = [double]"1.120183993975924e-28"
$ws.Range("J3").Value = [double]"3.412275579313158e-32"

# Row 4 (XY2)
$ws.Range("C4").Value = [double]"-2.857244234638895e-14"
$ws.Range("D4").Value = [double]"9.042249603192873e-14"
$ws.Range("E4").Value = [double]"-4.986831942807413e-16"
$ws.Range("F4").Value = [double]"1.578170273628664e-15"
$ws.Range("G4").Value = [double]"8.197870385363707e-29"
$ws.Range("H4").Value = [double]"2.497214124535363e-32"
$ws.Range("I4").Value = [double]"2.065801974244446e-28"
$ws.Range("J4").Value = [double]"6.292792671846319e-32"

# Row 5 (Y3)
$ws.Range("C5").Value = [double]"-7.160917981270024e-14"
$ws.Range("D5").Value = [double]"-1.822185867080097e-14"
$ws.Range("E5").Value = [double]"-1.249815962384275e-15"
$ws.Range("F5").Value = [double]"-3.180314296385544e-16"
$ws.Range("G5").Value = [double]"1.172710935399966e-28"
$ws.Range("H5").Value = [double]"3.572281792971952e-32"
$ws.Range("I5").Value = [double]"2.955143777208941e-28"
$ws.Range("J5").Value = [double]"9.001882725121354e-32"

# Row 6 (X2)
$ws.Range("C6").Value = [double]"1.588807724096984e-12"
$ws.Range("D6").Value = [double]"5.051750024577353e-11"
$ws.Range("E6").Value = [double]"2.772992596661002e-14"
$ws.Range("F6").Value = [double]"8.816967091657929e-13"
$ws.Range("G6").Value = [double]"3.555149200130355e-23"
$ws.Range("H6").Value = [double]"1.082960376300489e-26"
$ws.Range("I6").Value = [double]"8.958709873572875e-23"
$ws.Range("J6").Value = [double]"2.728979086305479e-26"

# Row 7 (XY)
$ws.Range("C7").Value = [double]"-1.393189389338127e-09"
$ws.Range("D7").Value = [double]"-2.365030851356797e-11"
$ws.Range("E7").Value = [double]"-2.43157419477995e-11"
$ws.Range("F7").Value = [double]"-4.127757526742071e-13"
$ws.Range("G7").Value = [double]"7.158917671638654e-23"
$ws.Range("H7").Value = [double]"2.180731029600038e-26"
$ws.Range("I7").Value = [double]"1.803993667175852e-22"
$ws.Range("J7").Value = [double]"5.495278962066702e-26"

# Row 8 (Y2)
$ws.Range("C8").Value = [double]"1.964924817432078e-11"
$ws.Range("D8").Value = [double]"-6.430609662071681e-10"
$ws.Range("E8").Value = [double]"3.429440761833822e-13"
$ws.Range("F8").Value = [double]"-1.122353115137108e-11"
$ws.Range("G8").Value = [double]"9.675012696024873e-23"
$ws.Range("H8").Value = [double]"2.947177403866744e-26"
$ws.Range("I8").Value = [double]"2.438030779795201e-22"
$ws.Range("J8").Value = [double]"7.426666455017913e-26"

# Row 9 (X)
$ws.Range("C9").Value = [double]"9.523697778139953e-06"
$ws.Range("D9").Value = [double]"-0.0005433124394329735"
$ws.Range("E9").Value = [double]"1.662198831934106e-07"
$ws.Range("F9").Value = [double]"-9.482590935147661e-06"
$ws.Range("G9").Value = [double]"2.888709817582726e-16"
$ws.Range("H9").Value = [double]"8.799513311445835e-20"
$ws.Range("I9").Value = [double]"7.279332514009924e-16"
$ws.Range("J9").Value = [double]"2.217411488187197e-19"

# Row 10 (Y)
$ws.Range("C10").Value = [double]"0.0005614393067380668"
$ws.Range("D10").Value = [double]"9.283324573376868e-06"
$ws.Range("E10").Value = [double]"9.798964452693649e-06"
$ws.Range("F10").Value = [double]"1.620245793367243e-07"
$ws.Range("G10").Value = [double]"4.826579192374805e-16"
$ws.Range("H10").Value = [double]"1.470260099977429e-19"
$ws.Range("I10").Value = [double]"1.216261828469087e-15"
$ws.Range("J10").Value = [double]"3.704945399733175e-19"

# Row 11 (ind)
$ws.Range("C11").Value = [double]"127.6475211903077"
$ws.Range("E11").Value = [double]"2.227869526780101"
$ws.Range("G11").Value = [double]"2.195491002298731e-10"
$ws.Range("H11").Value = [double]"6.687848042851741e-14"
$ws.Range("I11").Value = [double]"5.532472988450911e-10"
$ws.Range("J11").Value = [double]"1.685287646781577e-13"

Write-Host "Updated plate coefficients and variances."
